$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "add 2d act camera": give the village scene (row 2) its own camera
# offset/rotation instead of reusing the shared "0,8,7" / "45,180" pair.
$ws.Range("J2").Value = "0,4.2,5.5"
$ws.Range("K2").Value = "25,180"

$ws.Range("K2").Select()
